$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 0.1

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 0.2

$ws.Range("A4").Value = 15
$ws.Range("B4").Value = 0.3

$ws.Range("A5").Value = 20
$ws.Range("B5").Value = 0.4

$ws.Range("A6").Value = 25
$ws.Range("B6").Value = 0.5

$ws.Range("A7").Value = 30
$ws.Range("B7").Value = 0.6

$ws.Range("A8").Value = 35
$ws.Range("B8").Value = 0.7

$ws.Range("A9").Value = 40
$ws.Range("B9").Value = 0.8

$ws.Range("A10").Value = 50
$ws.Range("B10").Value = 1
